# Append a new row (row 7) to the equations sheet:
#   2^53 -> 9.007199e+15  at 1648538472436
#
# B7/C7 look numeric but must stay plain text (matching the rest of the
# sheet, which stores every column as text), so they are entered with a
# leading apostrophe to force text entry instead of numeric parsing.
# The explicit Style reset afterwards strips the "new cell inherits the
# column's style" / quote-prefix formatting Excel would otherwise leave
# behind, so row 7 ends up styled exactly like the other data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2^53"
$ws.Range("B7").Value = "'9.007199e+15"
$ws.Range("C7").Value = "'1648538472436"

$ws.Range("A7:C7").Style = "Normal"
